$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header labels in row 1:
#    columns A:J  "<Name>_old" -> "<Name>_FV2404"
#    column  K    "diff" stays unchanged
#    columns L:U  "<Name>_new" -> "<Name>_FV2410"
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $label = $cell.Value()
    $cell.Value = ($label -replace "_old$", "_FV2404")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $label = $cell.Value()
    $cell.Value = ($label -replace "_new$", "_FV2410")
}

# 2) Freeze the header row (row 1) so it stays visible while scrolling.
[void]$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true

# 3) Turn the data range A1:U54 into an Excel Table ("Table1") with an
#    autofilter, using the header row as column headers.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U54"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
